# 12.5.14 - Ceci Hours
#
# Add a new time-entry row to the "Work Database" table for Ceci:
#   Project : Reporting  - Git Hub
#   Sub-Task: Update
#   Month   : December
#   Year    : 2014
#   Hours   : 2
#   Person  : Ceci
#
# The new row grows the "Work Database" table/source range that feeds the
# "PivotTable2" pivot table on the "Summary" sheet, so the pivot is also
# refreshed to reflect Ceci's additional hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Database")

$row = 46
$ws.Cells.Item($row, 1).Value = "Reporting  - Git Hub"
$ws.Cells.Item($row, 2).Value = "Update"
$ws.Cells.Item($row, 3).Value = "December"
$ws.Cells.Item($row, 4).Value = 2014
$ws.Cells.Item($row, 5).Value = 2
$ws.Cells.Item($row, 6).Value = "Ceci"

# Refresh the pivot table / cache on the Summary sheet so it reflects the
# newly entered hours for Ceci.
$summary = $wb.Worksheets.Item("Summary")
$pivot = $summary.PivotTables(1)
$pivot.PivotCache().Refresh()
$pivot.RefreshTable()

# Leave the view roughly where the author left it: scrolled further down
# the "Work Database" sheet and back on that sheet as the active tab.
$summary.Activate()
$summary.Range("D45:D46").Select()

$ws.Activate()
$ws.Range("I44").Select()
